$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 33

$ws.Cells.Item($newRow, 1).Value = "2025/12/04 01:00"
$ws.Cells.Item($newRow, 2).Value = "-"
$ws.Cells.Item($newRow, 3).Value = "-"
$ws.Cells.Item($newRow, 4).Value = "-"
$ws.Cells.Item($newRow, 5).Value = "-"
$ws.Cells.Item($newRow, 6).Value = "-"
$ws.Cells.Item($newRow, 7).Value = "-"
